# Daily "symbol list" refresh for the crypto price sheet (GitHub Actions job).
# Only the Price column (D) moves for most rows; a handful of rows also pick
# up/lose the "Worstin24h"/"Bestin24h" badge text in column E, and three rows
# (41-43) got re-sorted by price so their Coin/Link/Price/Data cells rotate.
#
# Column D values are numeric-looking strings that must stay TEXT (the sheet
# stores prices as text, not numbers) -- prefixing with a leading apostrophe
# tells Excel to keep the literal text instead of coercing it to a Number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2  - BNB
$ws.Range("D2").Value = "'249.30"
# Row 3  - OKB
$ws.Range("D3").Value = "'22.70"
# Row 4  - HuobiToken
$ws.Range("D4").Value = "'5.298"
# Row 5  - Cronos
$ws.Range("D5").Value = "'0.05707"
# Row 6  - GateToken
$ws.Range("D6").Value = "'3.407"
# Row 7  - KuCoinToken
$ws.Range("D7").Value = "'6.335"
# Row 8  - MXToken
$ws.Range("D8").Value = "'0.8056"
# Row 9  - FTXToken
$ws.Range("D9").Value = "'0.8978"
# Row 10 - WazirX
$ws.Range("D10").Value = "'0.1425"
# Row 11 - MandalaExchangeToken
$ws.Range("D11").Value = "'0.07443"
# Row 12 - LiechtensteinCryptoassetsExchange
$ws.Range("D12").Value = "'0.03088"
# Row 13 - BitrueCoin
$ws.Range("D13").Value = "'0.03003"
# Row 14 - BitMartToken
$ws.Range("D14").Value = "'0.09408"
# Row 15 - MCDex
$ws.Range("D15").Value = "'3.867"
# Row 16 - BitForexToken
$ws.Range("D16").Value = "'0.001581"
# Row 17 - CoinExToken
$ws.Range("D17").Value = "'0.04794"
# Row 18 - UpBots
$ws.Range("D18").Value = "'0.01826"
# Row 19 - One (now flagged as the day's worst 24h performer)
$ws.Range("D19").Value = "'0.0005806"
$ws.Range("E19").Value = "18OneONEWorstin24h"
# Row 20 - TigerCash
$ws.Range("D20").Value = "'0.006417"
# Row 21 - HotbitToken
$ws.Range("D21").Value = "'0.004996"
# Row 22 - BitKan
$ws.Range("D22").Value = "'0.0009946"
# Row 23 - NitroEx
$ws.Range("D23").Value = "'0.0001498"
# Row 26 - BitpandaEcosystemToken
$ws.Range("D26").Value = "'0.3266"
# Row 27 - ProBitToken
$ws.Range("D27").Value = "'0.1354"
# Row 40 - IDEX
$ws.Range("D40").Value = "'0.03988"

# Rows 41-43 re-sorted by price: KickToken, BKEXToken, CEJI (in that order).
# Row 41 - was BKEXToken, now KickToken (also loses the Worstin24h badge)
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D41").Value = "'0.006783"
$ws.Range("E41").Value = "40KickTokenKICK"
# Row 42 - was CEJI, now BKEXToken
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").Value = "'0.1071"
$ws.Range("E42").Value = "41BKEXTokenBKK"
# Row 43 - was KickToken, now CEJI
$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D43").Value = "'0.002787"
$ws.Range("E43").Value = "42CEJICEJI"

# Row 44 - LocalTraders
$ws.Range("D44").Value = "'0.007680"
# Row 45 - CoinLion
$ws.Range("D45").Value = "'0.00005573"
# Row 47 - CoinbaseStockToken
$ws.Range("D47").Value = "'0.4987"
# Row 48 - BOLO
$ws.Range("D48").Value = "'0.2012"
